# Add partition() for small multiples
# Add a sheet to partition.xlsx for examples

$wb = $excel.ActiveWorkbook

# --- Preserve/record a selection on the "Pivot" sheet (matches the
# upstream edit where the Pivot sheet view moved to C36) without making
# it the active sheet.
$wsPivot = $wb.Worksheets.Item("Pivot")
$wsPivot.Range("C36").Select()

# --- Add the new "Small multiples" worksheet at the very end of the
# workbook (after "BELOW RIGHT border") and make it active.
$lastIndex = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($lastIndex))
$ws.Name = "Small multiples"

# --- Table 1: Postgraduate qualification | Bachelor's degree
$ws.Range("A1").Value = "Postgraduate qualification"
$ws.Range("D1").Value = "Bachelor's degree"

$ws.Range("A2").Value = "Sex"
$ws.Range("B2").Value = "Value"
$ws.Range("D2").Value = "Sex"
$ws.Range("E2").Value = "Value"

$ws.Range("A3").Value = "Female"
$ws.Range("B3").Value = 171000
$ws.Range("D3").Value = "Female"
$ws.Range("E3").Value = 275000

$ws.Range("A4").Value = "Male"
$ws.Range("B4").Value = 159000
$ws.Range("D4").Value = "Male"
$ws.Range("E4").Value = 200000

# --- Table 2: Diploma | Certificate
$ws.Range("A6").Value = "Diploma"
$ws.Range("D6").Value = "Certificate"

$ws.Range("A7").Value = "Sex"
$ws.Range("B7").Value = "Value"
$ws.Range("D7").Value = "Sex"
$ws.Range("E7").Value = "Value"

$ws.Range("A8").Value = "Female"
$ws.Range("B8").Value = 210000
$ws.Range("D8").Value = "Female"
$ws.Range("E8").Value = 732000

$ws.Range("A9").Value = "Male"
$ws.Range("B9").Value = 173000
$ws.Range("D9").Value = "Male"
$ws.Range("E9").Value = 807000

# --- Table 3: No Qualification
$ws.Range("A11").Value = "No Qualification"

$ws.Range("A12").Value = "Sex"
$ws.Range("B12").Value = "Value"

$ws.Range("A13").Value = "Female"
$ws.Range("B13").Value = 344000

$ws.Range("A14").Value = "Male"
$ws.Range("B14").Value = 287000

# --- Formatting: italic title rows, bold "Sex"/"Value" header rows
$ws.Range("A1:B1").Font.Italic = $true
$ws.Range("D1:E1").Font.Italic = $true
$ws.Range("A6:B6").Font.Italic = $true
$ws.Range("D6:E6").Font.Italic = $true
$ws.Range("A11:B11").Font.Italic = $true

$ws.Range("A2:B2").Font.Bold = $true
$ws.Range("D2:E2").Font.Bold = $true
$ws.Range("A7:B7").Font.Bold = $true
$ws.Range("D7:E7").Font.Bold = $true
$ws.Range("A12:B12").Font.Bold = $true

# --- Column widths (bestFit-like), approximated to the nearest width the
# host engine can represent.
$ws.Columns.Item(1).ColumnWidth = 13.0
$ws.Columns.Item(2).ColumnWidth = 6.0
$ws.Columns.Item(4).ColumnWidth = 21.33
$ws.Columns.Item(5).ColumnWidth = 6.0

# --- Selection state: land on E1 with the new sheet as the active tab.
$ws.Range("E1").Select()
